# "Ultimos ajustes para revision" - update the review-status columns on Hoja1
# with the notes from the latest meeting with Renzo (modales, encabezado, footer, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header y footer. (row 2): footer centering is now OK; drop the stray "OK Header" note ---
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "Pie ajustar centrado interior - OK"

# --- Contenido editorial. (row 6): gallery + gallery modal are now OK ---
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "Galeria - OK,  Modal Galería - OK"

# --- Posts. (row 8): element ordering issue resolved ---
$ws.Range("D8").Value = "OK"

# --- Home inicio. (row 10): notes from Carvajal's review of the featured article ---
$ws.Range("F10").Value = "Se realiza revision del articulo destacado -Carvajal: revisar estructura de #articulo-destacado, pues hay elementos dentro de un col-sx-6 y unos textos estan fuera."
$ws.Range("F10").WrapText = $true

# --- Row heights follow the new wrapped/trimmed text ---
$ws.Rows.Item(6).RowHeight = 28
$ws.Rows.Item(10).RowHeight = 42
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(26).RowHeight = 42

# --- Minor column width retuning ---
$ws.Columns.Item(1).ColumnWidth = 50.666666666666664
$ws.Columns.Item(2).ColumnWidth = 101.66666666666667
$ws.Columns.Item(5).ColumnWidth = 37.83
$ws.Columns.Item(6).ColumnWidth = 46.83

# --- Leave the cursor where the author left it ---
$ws.Range("B12").Select()
